# Updates the cryptos list data (prices, 1h volume %, and the two pairs of
# rows whose coins swapped rank) to match the refreshed scrape.
#
# Every D/E cell here holds plain text (prices formatted like "34.166.69"
# and percentages like "  +1.29%  "), so we force the Text number format
# before assigning the value (otherwise Excel would parse "0.548" etc. as
# a real number), then reset the style back to Normal so we don't leave a
# stray cell style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

$updates = @(
    @{ Cell = "D2";  Value = "34.161.06" }
    @{ Cell = "E2";  Value = "  +1.31%  " }

    @{ Cell = "D3";  Value = "1.787.18" }
    @{ Cell = "E3";  Value = "  +1.22%  " }

    @{ Cell = "D5";  Value = "226.57" }
    @{ Cell = "E5";  Value = "  +1.00%  " }

    @{ Cell = "D6";  Value = "0.548" }
    @{ Cell = "E6";  Value = "  +0.89%  " }

    @{ Cell = "E7";  Value = "  +0.05%  " }

    @{ Cell = "D8";  Value = "31.84" }
    @{ Cell = "E8";  Value = "  -0.26%  " }

    @{ Cell = "E9";  Value = "  +1.68%  " }

    @{ Cell = "E10"; Value = "  +0.49%  " }

    @{ Cell = "E11"; Value = "  +0.92%  " }

    @{ Cell = "D12"; Value = "2.045.79" }
    @{ Cell = "E12"; Value = "  +1.34%  " }

    # Rows 13 & 14 swapped coins (Chainlink <-> WrappedEther) with new data.
    @{ Cell = "B13"; Value = "WrappedEther" }
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" }
    @{ Cell = "D13"; Value = "1.796.99" }
    @{ Cell = "E13"; Value = "  +2.01%  " }

    @{ Cell = "B14"; Value = "Chainlink" }
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link" }
    @{ Cell = "D14"; Value = "11.04" }
    @{ Cell = "E14"; Value = "  -2.07%  " }

    @{ Cell = "E15"; Value = "  +2.34%  " }

    @{ Cell = "D16"; Value = "34.108.48" }
    @{ Cell = "E16"; Value = "  +1.22%  " }

    @{ Cell = "E17"; Value = "  +1.05%  " }

    @{ Cell = "D18"; Value = "68.25" }
    @{ Cell = "E18"; Value = "  +2.58%  " }

    @{ Cell = "D19"; Value = "246.87" }
    @{ Cell = "E19"; Value = "  +3.97%  " }

    @{ Cell = "E20"; Value = "  +0.81%  " }

    @{ Cell = "D23"; Value = "4.11" }
    @{ Cell = "E23"; Value = "  +1.57%  " }

    @{ Cell = "E24"; Value = "  +0.34%  " }

    @{ Cell = "D25"; Value = "161.15" }
    @{ Cell = "E25"; Value = "  +1.13%  " }

    @{ Cell = "E26"; Value = "  +2.48%  " }

    @{ Cell = "E27"; Value = "  +1.40%  " }

    @{ Cell = "E28"; Value = "  +1.42%  " }

    @{ Cell = "E29"; Value = "  +0.03%  " }

    @{ Cell = "E30"; Value = "  +0.78%  " }

    @{ Cell = "E31"; Value = "  +2.06%  " }

    @{ Cell = "E32"; Value = "  +2.96%  " }

    @{ Cell = "E33"; Value = "  +3.69%  " }

    @{ Cell = "E34"; Value = "  +0.99%  " }

    @{ Cell = "D35"; Value = "1.444.69" }
    @{ Cell = "E35"; Value = "  +4.72%  " }

    # Rows 36 & 37 swapped coins (RenderToken <-> ImmutableX) with new data.
    @{ Cell = "B36"; Value = "ImmutableX" }
    @{ Cell = "C36"; Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx" }
    @{ Cell = "D36"; Value = "0.654" }
    @{ Cell = "E36"; Value = "  +0.07%  " }

    @{ Cell = "B37"; Value = "RenderToken" }
    @{ Cell = "C37"; Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr" }
    @{ Cell = "D37"; Value = "2.45" }
    @{ Cell = "E37"; Value = "  +10.17%  " }

    @{ Cell = "E38"; Value = "  +4.04%  " }

    @{ Cell = "E39"; Value = "  +0.61%  " }

    @{ Cell = "D40"; Value = "80.68" }
    @{ Cell = "E40"; Value = "  +3.94%  " }

    @{ Cell = "E41"; Value = "  +0.88%  " }

    @{ Cell = "D42"; Value = "0.923" }
    @{ Cell = "E42"; Value = "  +1.93%  " }

    @{ Cell = "E43"; Value = "  +1.34%  " }

    @{ Cell = "D44"; Value = "13.58" }
    @{ Cell = "E44"; Value = "  -0.26%  " }

    @{ Cell = "D45"; Value = "6.07" }
    @{ Cell = "E45"; Value = "  +4.29%  " }

    @{ Cell = "E46"; Value = "  +1.75%  " }

    @{ Cell = "E47"; Value = "  -0.63%  " }

    @{ Cell = "E48"; Value = "  -2.17%  " }

    @{ Cell = "D49"; Value = "1.947.39" }
    @{ Cell = "E49"; Value = "  +1.57%  " }

    @{ Cell = "D50"; Value = "106.03" }
    @{ Cell = "E50"; Value = "  -1.47%  " }

    @{ Cell = "E51"; Value = "  +0.02%  " }
)

foreach ($u in $updates) {
    Set-TextValue $u.Cell $u.Value
}
